$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 33866.668
$ws.Range("I28").Value = 33866.668
$ws.Range("K28").Value = 33866.668
$ws.Range("M28").Value = -33381.668

$ws.Range("H76").Value = 71431910
$ws.Range("I76").Value = 3657.5715
$ws.Range("K76").Value = 3657.5715
$ws.Range("M76").Value = -3342.5715

$ws.Range("H79").Value = 71431910
$ws.Range("I79").Value = 3657.5715
$ws.Range("K79").Value = 3657.5715
$ws.Range("M79").Value = -2565.5715

$ws.Range("H113").Value = 2729.4
$ws.Range("I113").Value = 1882.6666
$ws.Range("J113").Value = 3999.5
$ws.Range("K113").Value = 1882.6666
$ws.Range("L113").Value = 3999.5
$ws.Range("M113").Value = 1371.3334
$ws.Range("N113").Value = -10507.5

$ws.Range("H135").Value = 281.2353
$ws.Range("I135").Value = 258.73334
$ws.Range("J135").Value = 450
$ws.Range("K135").Value = 2328.60006
$ws.Range("L135").Value = 4050
$ws.Range("M135").Value = 206.3999400000002
$ws.Range("N135").Value = -9120

$ws.Range("H137").Value = 1359.5
$ws.Range("I137").Value = 1353.7142
$ws.Range("J137").Value = 1400
$ws.Range("K137").Value = 4061.1426
$ws.Range("L137").Value = 4200
$ws.Range("M137").Value = -1511.1426
$ws.Range("N137").Value = -9300

$ws.Range("H138").Value = 1905.92
$ws.Range("I138").Value = 808.25
$ws.Range("J138").Value = 2115
$ws.Range("K138").Value = 2424.75
$ws.Range("L138").Value = 6345
$ws.Range("M138").Value = 2715.25
$ws.Range("N138").Value = -16625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9398.666999999999
$ws.Range("I2").Value = 971.6667
$ws.Range("K2").Value = 971.6667
$ws.Range("M2").Value = -858.6667

$ws.Range("H32").Value = 3144.8823
$ws.Range("I32").Value = 3230.3777
$ws.Range("J32").Value = 2503.6667
$ws.Range("K32").Value = 3230.3777
$ws.Range("L32").Value = 2503.6667
$ws.Range("M32").Value = -2943.3777
$ws.Range("N32").Value = -3077.6667

$ws.Range("H110").Value = 1869.0588
$ws.Range("I110").Value = 1554.3572
$ws.Range("J110").Value = 3337.6667
$ws.Range("K110").Value = 1554.3572
$ws.Range("L110").Value = 3337.6667
$ws.Range("M110").Value = 490.6428000000001
$ws.Range("N110").Value = -7427.6667

$ws.Range("H116").Value = 9398.666999999999
$ws.Range("I116").Value = 971.6667
$ws.Range("K116").Value = 971.6667
$ws.Range("M116").Value = 1322.3333

$ws.Range("H122").Value = 1135.3334
$ws.Range("I122").Value = 1135.3334
$ws.Range("K122").Value = 3406.0002
$ws.Range("M122").Value = -956.0001999999999

$ws.Range("H133").Value = 27970
$ws.Range("J133").Value = 27970
$ws.Range("L133").Value = 27970
$ws.Range("N133").Value = -33030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9398.666999999999
$ws.Range("I3").Value = 971.6667
$ws.Range("K3").Value = 971.6667
$ws.Range("M3").Value = -857.6667

$ws.Range("H80").Value = 986.6667
$ws.Range("I80").Value = 529.3333
$ws.Range("J80").Value = 1215.3334
$ws.Range("K80").Value = 529.3333
$ws.Range("L80").Value = 1215.3334
$ws.Range("M80").Value = 468.6667
$ws.Range("N80").Value = -3211.3334

$ws.Range("H83").Value = 986.6667
$ws.Range("I83").Value = 529.3333
$ws.Range("J83").Value = 1215.3334
$ws.Range("K83").Value = 2646.6665
$ws.Range("L83").Value = 6076.666999999999
$ws.Range("M83").Value = 2345.3335
$ws.Range("N83").Value = -16060.667

$ws.Range("H99").Value = 55557156
$ws.Range("I99").Value = 55557156
$ws.Range("K99").Value = 55557156
$ws.Range("M99").Value = -55555658

$ws.Range("H105").Value = 200001860
$ws.Range("I105").Value = 250001740
$ws.Range("J105").Value = 2309
$ws.Range("K105").Value = 250001740
$ws.Range("L105").Value = 2309
$ws.Range("M105").Value = -249999993
$ws.Range("N105").Value = -5803

$ws.Range("H107").Value = 1328.9678
$ws.Range("I107").Value = 935.58826
$ws.Range("J107").Value = 1806.6428
$ws.Range("K107").Value = 935.58826
$ws.Range("L107").Value = 1806.6428
$ws.Range("M107").Value = 984.41174
$ws.Range("N107").Value = -5646.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1569.4667
$ws.Range("I31").Value = 1467.2858
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1467.2858
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1172.2858
$ws.Range("N31").Value = -3590

$ws.Range("H34").Value = 1569.4667
$ws.Range("I34").Value = 1467.2858
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1467.2858
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1265.2858
$ws.Range("N34").Value = -3404

$ws.Range("H58").Value = 624.7826
$ws.Range("I58").Value = 607.0909
$ws.Range("J58").Value = 1014
$ws.Range("K58").Value = 607.0909
$ws.Range("L58").Value = 1014
$ws.Range("M58").Value = -404.0909
$ws.Range("N58").Value = -1420

$ws.Range("H99").Value = 1445.7646
$ws.Range("I99").Value = 1427.2142
$ws.Range("K99").Value = 1427.2142
$ws.Range("M99").Value = 70.78580000000011

$ws.Range("H126").Value = 1445.7646
$ws.Range("I126").Value = 1427.2142
$ws.Range("K126").Value = 4281.642599999999
$ws.Range("M126").Value = -1811.642599999999

$ws.Range("H132").Value = 4720.6333
$ws.Range("I132").Value = 4800.7407
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 14402.2221
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -11872.2221
$ws.Range("N132").Value = -17059.0001

$ws.Range("H134").Value = 27780132
$ws.Range("I134").Value = 41668972
$ws.Range("J134").Value = 2450
$ws.Range("K134").Value = 125006916
$ws.Range("L134").Value = 7350
$ws.Range("M134").Value = -125004381
$ws.Range("N134").Value = -12420

$ws.Range("H136").Value = 624.7826
$ws.Range("I136").Value = 607.0909
$ws.Range("J136").Value = 1014
$ws.Range("K136").Value = 1821.2727
$ws.Range("L136").Value = 3042
$ws.Range("M136").Value = 728.7273
$ws.Range("N136").Value = -8142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 856.5
$ws.Range("I122").Value = 680.7143
$ws.Range("J122").Value = 1266.6666
$ws.Range("K122").Value = 6126.428699999999
$ws.Range("L122").Value = 11399.9994
$ws.Range("M122").Value = -3676.428699999999
$ws.Range("N122").Value = -16299.9994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1966

$ws.Range("H71").Value = 1966

$ws.Range("H93").Value = 548.3077
$ws.Range("I93").Value = 502.1111
$ws.Range("J93").Value = 652.25
$ws.Range("K93").Value = 502.1111
$ws.Range("L93").Value = 652.25
$ws.Range("M93").Value = 745.8888999999999
$ws.Range("N93").Value = -3148.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 297.7143
$ws.Range("I107").Value = 313.1579
$ws.Range("J107").Value = 265.1111
$ws.Range("K107").Value = 939.4737
$ws.Range("L107").Value = 795.3333
$ws.Range("M107").Value = 980.5263
$ws.Range("N107").Value = -4635.3333

$ws.Range("H122").Value = 7030502.5
$ws.Range("I122").Value = 8128496.5
$ws.Range("K122").Value = 24385489.5
$ws.Range("M122").Value = -24383039.5

$ws.Range("H132").Value = 3293.28
$ws.Range("I132").Value = 3270.6875
$ws.Range("J132").Value = 3333.4443
$ws.Range("K132").Value = 9812.0625
$ws.Range("L132").Value = 10000.3329
$ws.Range("M132").Value = -7282.0625
$ws.Range("N132").Value = -15060.3329
